# Sets ELC-technology availability for the new regions
# (Availability sheet: shift technology rows, add AllRegions wildcard row,
#  clear the now-unused SUPH2* rows; also tidy up sheet scroll/selection state)

$wb = $excel.ActiveWorkbook
$wsAvail = $wb.Worksheets.Item("Availability")
$wsElc = $wb.Worksheets.Item("ELC_Constraint")

# --- Availability sheet: row 5 becomes the new "AllRegions" wildcard row ---
$wsAvail.Range("C5").Value = 1
$wsAvail.Range("F5").Value = 0
$wsAvail.Range("G5").Value = 0
$wsAvail.Range("H5").Value = 0
$wsAvail.Range("I5").Value = 0
$wsAvail.Range("J5").Value = "*"

# --- Row 6: was ERWINELCWIN5N (style 11) -> becomes ERWINELCWIN3N (no explicit style) ---
$wsAvail.Range("F6:J6").ClearFormats()
$wsAvail.Range("F6").Value = 1
$wsAvail.Range("G6").Value = 1
$wsAvail.Range("H6").Value = 1
$wsAvail.Range("I6").Value = 1
$wsAvail.Range("J6").Value = "ERWINELCWIN3N"

# --- Row 7: was SUPH2ALKC1N (style 42) -> becomes ERWINELCWIN5N (style 11) ---
$wsAvail.Range("F7").Value = 1
$wsAvail.Range("G7").Value = 1
$wsAvail.Range("H7").Value = 1
$wsAvail.Range("I7").Value = 1
$wsAvail.Range("J7").Value = "ERWINELCWIN5N"
$wsAvail.Range("J7").ClearFormats()
$wsAvail.Range("J7").Style = $wsAvail.Range("F7").Style

# --- Rows 8-11: old SUPH2* technology rows no longer apply -> clear values, keep formatting ---
$wsAvail.Range("F8:J11").ClearContents()

# --- View/selection bookkeeping ---
$wsElc.Activate()
$excel.ActiveWindow.ScrollColumn = 1

$wsAvail.Activate()
$wsAvail.Range("I8").Select()
